$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ghi cong")
$ws.Activate()
$app = $ws.Application
$win = $app.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("H13").Select()
